$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.2594890654562015
$ws.Range("D2").Value = -0.2934520523824585
$ws.Range("E2").Value = -0.2494377539355348
$ws.Range("F2").Value = -0.03617483367244573
$ws.Range("G2").Value = 0.09003991013025806
$ws.Range("H2").Value = 0.1458347287376139
$ws.Range("I2").Value = -0.1399801602867078
$ws.Range("J2").Value = 0.01939071256390057
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -0.2926153968279041
$ws.Range("M2").Value = -0.05495139749410272
$ws.Range("N2").Value = 0.3388036330279748
# Row 3
$ws.Range("C3").Value = -0.4185302256558379
$ws.Range("D3").Value = -0.2234835412280841
$ws.Range("E3").Value = -0.2072920714574322
$ws.Range("F3").Value = -0.07166315892348649
$ws.Range("G3").Value = 0.2788672574637895
$ws.Range("H3").Value = 0.1973950443594833
$ws.Range("I3").Value = -0.7127793382414593
$ws.Range("J3").Value = 0.6341125469699391
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = -1.608942561718531
$ws.Range("M3").Value = 0.7896345490914247
$ws.Range("N3").Value = 0.8703042820675433
# Row 4
$ws.Range("C4").Value = -1.22006804823917
$ws.Range("D4").Value = 0.010570696466361
$ws.Range("E4").Value = -2.050652521105039
$ws.Range("F4").Value = 0.1750662512221364
$ws.Range("G4").Value = 1.783504520762116
$ws.Range("H4").Value = 0.06828100991706487
$ws.Range("I4").Value = -1.090443284667019
$ws.Range("J4").Value = 0.3246661265630818
$ws.Range("K4").Value = -0
$ws.Range("L4").Value = -1.001199160376046
$ws.Range("M4").Value = 0.7158355882672631
$ws.Range("N4").Value = -0.8639221881068608
# Row 5
$ws.Range("C5").Value = -0.9789051338303025
$ws.Range("D5").Value = 0.2832876626337821
$ws.Range("E5").Value = -2.140514980176671
$ws.Range("F5").Value = 0.411452656807904
$ws.Range("G5").Value = 1.448833769525207
$ws.Range("H5").Value = -0.5686922471844291
$ws.Range("I5").Value = -1.256765309321006
$ws.Range("J5").Value = 0.6025108545893475
$ws.Range("K5").Value = -0
$ws.Range("L5").Value = -1.746787450546593
$ws.Range("M5").Value = 1.301702462471367
$ws.Range("N5").Value = -0.2285659650834222
# Row 6
$ws.Range("C6").Value = -1.947080427264786
$ws.Range("D6").Value = 1.097773124942438
$ws.Range("E6").Value = -2.678569806967828
$ws.Range("F6").Value = 1.482551069300036
$ws.Range("G6").Value = 2.134907188550496
$ws.Range("H6").Value = -1.548182035994453
$ws.Range("I6").Value = -3.824684060090842
$ws.Range("J6").Value = 0.8612236576940775
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = -5.650182057401931
$ws.Range("M6").Value = 4.202965842269197
$ws.Range("N6").Value = -1.653149652953792
# Row 7
$ws.Range("C7").Value = -0.3444362924635327
$ws.Range("D7").Value = 0.1681585715117228
$ws.Range("E7").Value = -1.298318724003343
$ws.Range("F7").Value = 0.2359069376360411
$ws.Range("G7").Value = 0.3364880421999878
$ws.Range("H7").Value = -0.249257114077595
$ws.Range("I7").Value = -0.6730524596865569
$ws.Range("J7").Value = 0.1520416095900123
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = -0.919396568622477
$ws.Range("M7").Value = 0.7230866621558039
$ws.Range("N7").Value = -0.3029257110927108
# Row 8
$ws.Range("C8").Value = 0.03468264945443265
$ws.Range("D8").Value = 0.02402384500745207
$ws.Range("E8").Value = -0.3495953757225437
$ws.Range("F8").Value = 0.02183415705656238
$ws.Range("G8").Value = -0.006850299401198374
$ws.Range("H8").Value = 0.004684701668394968
$ws.Range("I8").Value = -0.02011187388761828
$ws.Range("J8").Value = -0.0242115835286234
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = -0.1063643457903465
$ws.Range("M8").Value = 0.02699550074987656
# Row 9
$ws.Range("C9").Value = 0.2253262209110191
$ws.Range("D9").Value = 0.1136597030652758
$ws.Range("E9").Value = -0.4152909437775469
$ws.Range("F9").Value = 0.02095401642685323
$ws.Range("G9").Value = 0.05446279599041853
$ws.Range("H9").Value = -0.04212149311311147
$ws.Range("I9").Value = -0.1534334619483862
$ws.Range("J9").Value = 0.04028556482922666
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = -0.9621218509655438
$ws.Range("M9").Value = 0.6915806363956095
# Row 10
$ws.Range("C10").Value = 0.1122520395495311
$ws.Range("D10").Value = 0.137652844838152
$ws.Range("E10").Value = -1.644385879891446
$ws.Range("F10").Value = -0.06439117150823856
$ws.Range("G10").Value = -0.2772284572608367
$ws.Range("H10").Value = -0.03416341761919085
$ws.Range("I10").Value = 0.03198665965268126
$ws.Range("J10").Value = -0.1392669483330781
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = -1.788884105880865
$ws.Range("M10").Value = 0.9954024011598053
# Row 11
$ws.Range("C11").Value = 0.2389878252173965
$ws.Range("D11").Value = 0.1636148486915047
$ws.Range("E11").Value = -1.651636982109398
$ws.Range("F11").Value = -0.05078685021685899
$ws.Range("G11").Value = -0.1382479055080232
$ws.Range("H11").Value = -0.1118639946266376
$ws.Range("I11").Value = -0.02071681530078861
$ws.Range("J11").Value = -0.04300018881633662
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = -2.605765297647407
$ws.Range("M11").Value = 1.407403986358935
# Row 12
$ws.Range("C12").Value = 0.8990615323368752
$ws.Range("D12").Value = 0.4632329975338321
$ws.Range("E12").Value = -1.145093050071671
$ws.Range("F12").Value = 0.4245736807894576
$ws.Range("G12").Value = 0.4386466909671858
$ws.Range("H12").Value = -0.2937527160706956
$ws.Range("I12").Value = -0.5978546100523756
$ws.Range("J12").Value = 0.8805452787703727
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = -5.768774203058401
$ws.Range("M12").Value = 1.652283729358433
# Row 13
$ws.Range("C13").Value = 0.1596310907246765
$ws.Range("D13").Value = 0.07834268303374968
$ws.Range("E13").Value = -0.3051681095420896
$ws.Range("F13").Value = 0.07691857335908939
$ws.Range("G13").Value = 0.07862533823311679
$ws.Range("H13").Value = -0.05125319195938525
$ws.Range("I13").Value = -0.1040808282978039
$ws.Range("J13").Value = 0.1550415574641953
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = -0.9351569015473625
$ws.Range("M13").Value = 0.3111619189530186
# Row 14
$ws.Range("C14").Value = -0.2941717149106342
$ws.Range("D14").Value = -0.3174758973899106
$ws.Range("E14").Value = 0.1001576217870089
$ws.Range("F14").Value = -0.05800899072900811
$ws.Range("G14").Value = 0.09689020953145643
$ws.Range("H14").Value = 0.141150027069219
$ws.Range("I14").Value = -0.1198682863990895
$ws.Range("J14").Value = 0.04360229609252397
$ws.Range("L14").Value = -0.1862510510375576
$ws.Range("M14").Value = -0.08194689824397927
$ws.Range("N14").Value = 0.3388036330279748
# Row 15
$ws.Range("C15").Value = -0.6438564465668569
$ws.Range("D15").Value = -0.3371432442933599
$ws.Range("E15").Value = 0.2079988723201147
$ws.Range("F15").Value = -0.09261717535033971
$ws.Range("G15").Value = 0.224404461473371
$ws.Range("H15").Value = 0.2395165374725948
$ws.Range("I15").Value = -0.5593458762930731
$ws.Range("J15").Value = 0.5938269821407125
$ws.Range("L15").Value = -0.6468207107529873
$ws.Range("M15").Value = 0.0980539126958152
$ws.Range("N15").Value = 0.8703042820675433
# Row 16
$ws.Range("C16").Value = -1.332320087788701
$ws.Range("D16").Value = -0.127082148371791
$ws.Range("E16").Value = -0.4062666412135938
$ws.Range("F16").Value = 0.239457422730375
$ws.Range("G16").Value = 2.060732978022953
$ws.Range("H16").Value = 0.1024444275362557
$ws.Range("I16").Value = -1.1224299443197
$ws.Range("J16").Value = 0.46393307489616
$ws.Range("K16").Value = -0
$ws.Range("L16").Value = 0.7876849455048192
$ws.Range("M16").Value = -0.2795668128925421
$ws.Range("N16").Value = -0.8639221881068608
# Row 17
$ws.Range("C17").Value = -1.217892959047699
$ws.Range("D17").Value = 0.1196728139422774
$ws.Range("E17").Value = -0.4888779980672731
$ws.Range("F17").Value = 0.462239507024763
$ws.Range("G17").Value = 1.587081675033231
$ws.Range("H17").Value = -0.4568282525577915
$ws.Range("I17").Value = -1.236048494020218
$ws.Range("J17").Value = 0.6455110434056841
$ws.Range("K17").Value = -0
$ws.Range("L17").Value = 0.8589778471008143
$ws.Range("M17").Value = -0.105701523887568
$ws.Range("N17").Value = -0.2285659650834222
# Row 18
$ws.Range("C18").Value = -2.846141959601661
$ws.Range("D18").Value = 0.6345401274086058
$ws.Range("E18").Value = -1.533476756896157
$ws.Range("F18").Value = 1.057977388510578
$ws.Range("G18").Value = 1.69626049758331
$ws.Range("H18").Value = -1.254429319923758
$ws.Range("I18").Value = -3.226829450038467
$ws.Range("J18").Value = -0.0193216210762952
$ws.Range("L18").Value = 0.1185921456564696
$ws.Range("M18").Value = 2.550682112910764
$ws.Range("N18").Value = -1.653149652953792
# Row 19
$ws.Range("C19").Value = -0.5040673831882092
$ws.Range("D19").Value = 0.08981588847797307
$ws.Range("E19").Value = -0.9931506144612535
$ws.Range("F19").Value = 0.1589883642769518
$ws.Range("G19").Value = 0.2578627039668711
$ws.Range("H19").Value = -0.1980039221182097
$ws.Range("I19").Value = -0.568971631388753
$ws.Range("J19").Value = -0.002999947874183007
$ws.Range("L19").Value = 0.01576033292488555
$ws.Range("M19").Value = 0.4119247432027853
$ws.Range("N19").Value = -0.3029257110927108
